$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily UF/IVP rows for 10-08-2021 .. 09-09-2021, appended after the
# existing data (which ends at row 222 / 09-08-2021).
$newRows = @(
    @{ Row = 223; Date = "10-08-2021"; UF = 29773.93; IVP = 31000.74 }
    @{ Row = 224; Date = "11-08-2021"; UF = 29781.59; IVP = 31004.32 }
    @{ Row = 225; Date = "12-08-2021"; UF = 29789.24; IVP = 31007.9 }
    @{ Row = 226; Date = "13-08-2021"; UF = 29796.9; IVP = 31011.48 }
    @{ Row = 227; Date = "14-08-2021"; UF = 29804.56; IVP = 31015.06 }
    @{ Row = 228; Date = "15-08-2021"; UF = 29812.22; IVP = 31018.64 }
    @{ Row = 229; Date = "16-08-2021"; UF = 29819.89; IVP = 31022.22 }
    @{ Row = 230; Date = "17-08-2021"; UF = 29827.55; IVP = 31025.8 }
    @{ Row = 231; Date = "18-08-2021"; UF = 29835.22; IVP = 31029.39 }
    @{ Row = 232; Date = "19-08-2021"; UF = 29842.89; IVP = 31032.97 }
    @{ Row = 233; Date = "20-08-2021"; UF = 29850.56; IVP = 31036.55 }
    @{ Row = 234; Date = "21-08-2021"; UF = 29858.23; IVP = 31040.14 }
    @{ Row = 235; Date = "22-08-2021"; UF = 29865.91; IVP = 31043.72 }
    @{ Row = 236; Date = "23-08-2021"; UF = 29873.59; IVP = 31047.31 }
    @{ Row = 237; Date = "24-08-2021"; UF = 29881.27; IVP = 31050.89 }
    @{ Row = 238; Date = "25-08-2021"; UF = 29888.95; IVP = 31054.48 }
    @{ Row = 239; Date = "26-08-2021"; UF = 29896.63; IVP = 31058.06 }
    @{ Row = 240; Date = "27-08-2021"; UF = 29904.32; IVP = 31061.65 }
    @{ Row = 241; Date = "28-08-2021"; UF = 29912.01; IVP = 31065.23 }
    @{ Row = 242; Date = "29-08-2021"; UF = 29919.7; IVP = 31068.82 }
    @{ Row = 243; Date = "30-08-2021"; UF = 29927.39; IVP = 31072.41 }
    @{ Row = 244; Date = "31-08-2021"; UF = 29935.08; IVP = 31076 }
    @{ Row = 245; Date = "01-09-2021"; UF = 29942.78; IVP = 31079.59 }
    @{ Row = 246; Date = "02-09-2021"; UF = 29950.47; IVP = 31083.17 }
    @{ Row = 247; Date = "03-09-2021"; UF = 29958.17; IVP = 31086.76 }
    @{ Row = 248; Date = "04-09-2021"; UF = 29965.87; IVP = 31090.35 }
    @{ Row = 249; Date = "05-09-2021"; UF = 29973.58; IVP = 31093.94 }
    @{ Row = 250; Date = "06-09-2021"; UF = 29981.28; IVP = 31097.53 }
    @{ Row = 251; Date = "07-09-2021"; UF = 29988.99; IVP = 31101.12 }
    @{ Row = 252; Date = "08-09-2021"; UF = 29996.7; IVP = 31104.71 }
    @{ Row = 253; Date = "09-09-2021"; UF = 30004.41; IVP = 31108.31 }
)

foreach ($r in $newRows) {
    $dateCell = $ws.Cells.Item($r.Row, 1)
    # Force the "dd-mm-yyyy" literal to be stored as text (matching the
    # other rows in column A) instead of being auto-converted into a date
    # serial number; ClearFormats afterwards removes the temporary "@"
    # number format again so the cell keeps the workbook default style.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r.Date
    $dateCell.ClearFormats()

    $ws.Cells.Item($r.Row, 2).Value = $r.UF
    $ws.Cells.Item($r.Row, 3).Value = $r.IVP
}
